$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" date column (C) for rows 2-5 from 45175 to 45183
$ws.Range("C2").Value = 45183
$ws.Range("C3").Value = 45183
$ws.Range("C4").Value = 45183
$ws.Range("C5").Value = 45183
